$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152, pushing existing rows 152-202 down to 153-203.
$ws.Rows.Item(152).Insert()

# Populate the new row 152 with the new data record.
$ws.Range("A152").Value = 3
$ws.Range("B152").Value = "Femacal de La Calera"
$ws.Range("C152").Value = "Coquimbo"
$ws.Range("D152").Value = 44900
$ws.Range("E152").Value = 5
$ws.Range("F152").Value = 100112030
$ws.Range("G152").Value = "Poroto granado"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 73
$ws.Range("K152").Value = 38000
$ws.Range("L152").Value = 39000
$ws.Range("M152").Value = 38479
$ws.Range("N152").Value = "$/saco 25 kilos"
$ws.Range("O152").Value = "Provincia de Limarí"
$ws.Range("P152").Value = 1539
$ws.Range("Q152").Value = 25
$ws.Range("R152").Value = "Hortaliza"
